$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 7339
$ws.Range("N10").Value = -1647
$ws.Range("L10").Value = 1061
$ws.Range("J10").Value = 1061
$ws.Range("N18").Value = -15566.667
$ws.Range("J18").Value = 14998.667
$ws.Range("H18").Value = 14998.667
$ws.Range("L18").Value = 14998.667
$ws.Range("J62").Value = 22679.143
$ws.Range("M62").Value = -13540211
$ws.Range("I62").Value = 13540835
$ws.Range("K62").Value = 13540835
$ws.Range("L62").Value = 22679.143
$ws.Range("H62").Value = 8809481
$ws.Range("N62").Value = -23927.143
$ws.Range("J65").Value = 22679.143
$ws.Range("L65").Value = 113395.715
$ws.Range("I65").Value = 13540835
$ws.Range("H65").Value = 8809481
$ws.Range("N65").Value = -119635.715
$ws.Range("M65").Value = -67701055
$ws.Range("K65").Value = 67704175
$ws.Range("H86").Value = 5816.3335
$ws.Range("M86").Value = -5056.6
$ws.Range("I86").Value = 6179.6
$ws.Range("K86").Value = 6179.6
$ws.Range("I89").Value = 6179.6
$ws.Range("M89").Value = -25282
$ws.Range("K89").Value = 30898
$ws.Range("H89").Value = 5816.3335
$ws.Range("H107").Value = 1070.7778
$ws.Range("I107").Value = 1094.1428
$ws.Range("M107").Value = 825.8571999999999
$ws.Range("K107").Value = 1094.1428
$ws.Range("K113").Value = 8247.5
$ws.Range("J113").Value = 4992.4287
$ws.Range("M113").Value = -4993.5
$ws.Range("L113").Value = 4992.4287
$ws.Range("I113").Value = 8247.5
$ws.Range("H113").Value = 6176.091
$ws.Range("N113").Value = -11500.4287
$ws.Range("L129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("H129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H137").Value = 2409565.2
$ws.Range("M137").Value = -9193999.5
$ws.Range("J137").Value = 4410.1113
$ws.Range("L137").Value = 13230.3339
$ws.Range("N137").Value = -18330.3339
$ws.Range("I137").Value = 3065516.5
$ws.Range("K137").Value = 9196549.5
$ws.Range("I138").Value = 3828.3333
$ws.Range("L138").Value = 13320
$ws.Range("J138").Value = 4440
$ws.Range("M138").Value = -6344.999899999999
$ws.Range("K138").Value = 11484.9999
$ws.Range("N138").Value = -23600
$ws.Range("H138").Value = 3915.7144
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 1179.625
$ws.Range("N32").Value = -1753.625
$ws.Range("L32").Value = 1179.625
$ws.Range("H32").Value = 2241.5696
$ws.Range("M61").Value = -4252
$ws.Range("N61").Value = -10694
$ws.Range("H61").Value = 5991.8945
$ws.Range("K61").Value = 4464
$ws.Range("I61").Value = 4464
$ws.Range("J61").Value = 10270
$ws.Range("L61").Value = 10270
$ws.Range("I63").Value = 1893.9375
$ws.Range("H63").Value = 1961.2778
$ws.Range("K63").Value = 1893.9375
$ws.Range("M63").Value = -1207.9375
$ws.Range("H66").Value = 1961.2778
$ws.Range("I66").Value = 1893.9375
$ws.Range("M66").Value = -6037.6875
$ws.Range("K66").Value = 9469.6875
$ws.Range("H74").Value = 56452.582
$ws.Range("L74").Value = 4309.8
$ws.Range("N74").Value = -6057.8
$ws.Range("J74").Value = 4309.8
$ws.Range("H77").Value = 56452.582
$ws.Range("L77").Value = 21549
$ws.Range("N77").Value = -30285
$ws.Range("J77").Value = 4309.8
$ws.Range("J125").Value = 51715
$ws.Range("L125").Value = 51715
$ws.Range("N125").Value = -61555
$ws.Range("H125").Value = 51715
$ws.Range("I132").Value = 1826
$ws.Range("K132").Value = 5478
$ws.Range("M132").Value = -2948
$ws.Range("H132").Value = 2109.2
$ws.Range("K136").Value = 13392
$ws.Range("H136").Value = 5991.8945
$ws.Range("L136").Value = 30810
$ws.Range("N136").Value = -35910
$ws.Range("J136").Value = 10270
$ws.Range("I136").Value = 4464
$ws.Range("M136").Value = -10842
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M20").Value = -4841.75
$ws.Range("K20").Value = 5088.75
$ws.Range("H20").Value = 4966.6
$ws.Range("I20").Value = 5088.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("M22").Value = 148.75
$ws.Range("K22").Value = 201.25
$ws.Range("I22").Value = 201.25
$ws.Range("H22").Value = 290
$ws.Range("N68").Value = -44498
$ws.Range("J68").Value = 43000
$ws.Range("L68").Value = 43000
$ws.Range("H68").Value = 43000
$ws.Range("L71").Value = 129000
$ws.Range("J71").Value = 43000
$ws.Range("H71").Value = 43000
$ws.Range("N71").Value = -136488
$ws.Range("I122").Value = 732.875
$ws.Range("M122").Value = 251.375
$ws.Range("L122").Value = 2368.8
$ws.Range("J122").Value = 789.6
$ws.Range("K122").Value = 2198.625
$ws.Range("H122").Value = 754.6923
$ws.Range("N122").Value = -7268.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I2").Value = 47.92857
$ws.Range("M2").Value = -174.57142
$ws.Range("K2").Value = 287.57142
$ws.Range("H2").Value = 90.652176
$ws.Range("L6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("H6").Value = 636.75
$ws.Range("N6").ClearContents()
$ws.Range("H7").Value = 333413.66
$ws.Range("M7").Value = -1499946.5
$ws.Range("J7").Value = 202
$ws.Range("I7").Value = 500019.5
$ws.Range("L7").Value = 606
$ws.Range("K7").Value = 1500058.5
$ws.Range("N7").Value = -830
$ws.Range("M11").Value = -5883.4
$ws.Range("K11").Value = 6023.4
$ws.Range("H11").Value = 2007.8
$ws.Range("I11").Value = 2007.8
$ws.Range("H12").Value = 4602.875
$ws.Range("J12").Value = 6450.091
$ws.Range("L12").Value = 19350.273
$ws.Range("N12").Value = -19696.273
$ws.Range("L13").Value = 9000
$ws.Range("N13").Value = -9336
$ws.Range("I13").Value = 1640.5
$ws.Range("M13").Value = -4753.5
$ws.Range("J13").Value = 3000
$ws.Range("H13").Value = 2093.6667
$ws.Range("K13").Value = 4921.5
$ws.Range("J15").Value = 674
$ws.Range("L15").Value = 2022
$ws.Range("H15").Value = 739
$ws.Range("N15").Value = -2302
$ws.Range("I17").Value = 143.75
$ws.Range("H17").Value = 143.75
$ws.Range("M17").Value = -262.25
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("K17").Value = 431.25
$ws.Range("N17").ClearContents()
$ws.Range("N19").Value = -6350.625
$ws.Range("J19").Value = 2000.875
$ws.Range("L19").Value = 6002.625
$ws.Range("H19").Value = 2000.875
$ws.Range("I23").Value = 861.6667
$ws.Range("L23").Value = 3920.1819
$ws.Range("H23").Value = 1211.3572
$ws.Range("N23").Value = -4390.1819
$ws.Range("J23").Value = 1306.7273
$ws.Range("M23").Value = -2350.0001
$ws.Range("K23").Value = 2585.0001
$ws.Range("N131").Value = -15089.7789
$ws.Range("J131").Value = 1669.9263
$ws.Range("L131").Value = 5009.7789
$ws.Range("H131").Value = 1663.3093
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J2").Value = 161
$ws.Range("H2").Value = 151.21428
$ws.Range("L2").Value = 161
$ws.Range("N2").Value = -387
$ws.Range("J18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("K70").Value = 6895.5625
$ws.Range("M70").Value = -6625.5625
$ws.Range("H70").Value = 7820.4
$ws.Range("I70").Value = 6895.5625
$ws.Range("K73").Value = 6895.5625
$ws.Range("M73").Value = -5959.5625
$ws.Range("I73").Value = 6895.5625
$ws.Range("H73").Value = 7820.4
$ws.Range("J102").Value = 3289.5
$ws.Range("N102").Value = -6533.5
$ws.Range("H102").Value = 2870.88
$ws.Range("L102").Value = 3289.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L21").Value = 13500
$ws.Range("H21").Value = 13500
$ws.Range("J21").Value = 13500
$ws.Range("N21").Value = -13848
$ws.Range("I53").Value = 15750
$ws.Range("L53").Value = 0
$ws.Range("H53").Value = 15750
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 15750
$ws.Range("M53").Value = -15232
$ws.Range("N53").ClearContents()
$ws.Range("J55").Value = 12839.833
$ws.Range("I55").Value = 793.2778
$ws.Range("N55").Value = -13185.833
$ws.Range("H55").Value = 5611.9
$ws.Range("K55").Value = 793.2778
$ws.Range("L55").Value = 12839.833
$ws.Range("M55").Value = -620.2778
$ws.Range("K136").Value = 8549.625
$ws.Range("H136").Value = 2873.111
$ws.Range("L136").Value = 8720.7276
$ws.Range("N136").Value = -13820.7276
$ws.Range("J136").Value = 2906.9092
$ws.Range("I136").Value = 2849.875
$ws.Range("M136").Value = -5999.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I23").Value = 832.5
$ws.Range("L23").Value = 3000
$ws.Range("H23").Value = 1266
$ws.Range("N23").Value = -3458
$ws.Range("J23").Value = 3000
$ws.Range("M23").Value = -603.5
$ws.Range("K23").Value = 832.5
$ws.Range("J33").Value = 6950
$ws.Range("N33").Value = -7450
$ws.Range("L33").Value = 6950
$ws.Range("H33").Value = 6950
$ws.Range("J36").Value = 6950
$ws.Range("H36").Value = 6950
$ws.Range("N36").Value = -7450
$ws.Range("L36").Value = 6950
$ws.Range("H107").Value = 666.3
$ws.Range("I107").Value = 558.92
$ws.Range("M107").Value = 243.2400000000002
$ws.Range("K107").Value = 1676.76
$ws.Range("I132").Value = 2333
$ws.Range("K132").Value = 6999
$ws.Range("M132").Value = -4469
$ws.Range("H132").Value = 2333
